# Fruta / hortaliza, semanal
# A new weekly price observation for "Alcachofa" (Macroferia Regional de Talca)
# is inserted at the top of the existing series (row 65), pushing the prior
# rows (65-135) down by one (66-136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 65, shifting rows 65..135 down to 66..136.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with this week's observation.
$ws.Range("A65").Value = 5
$ws.Range("B65").Value = "Macroferia Regional de Talca"
$ws.Range("C65").Value = "Maule"
$ws.Range("D65").Value = "2023-08-21"
$ws.Range("E65").Value = 7
$ws.Range("F65").Value = 100112013
$ws.Range("G65").Value = "Alcachofa"
$ws.Range("H65").Value = "Madrigal"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 400
$ws.Range("K65").Value = 13000
$ws.Range("L65").Value = 13000
$ws.Range("M65").Value = 13000
$ws.Range("N65").Value = "$/caja 40 unidades"
$ws.Range("O65").Value = "Provincia del Elquí"
$ws.Range("P65").Value = 325
$ws.Range("Q65").Value = 40
$ws.Range("R65").Value = "Hortaliza"
